# Apply updated symbol/price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'318.67"
$ws.Range("E2").Value = "'3.89%"
$ws.Range("D3").Value = "'39.75"
$ws.Range("E3").Value = "'2.27%"
$ws.Range("D4").Value = "'5.142"
$ws.Range("E4").Value = "'0.89%"
$ws.Range("D5").Value = "'0.08221"
$ws.Range("E5").Value = "'1.84%"
$ws.Range("D6").Value = "'2.056"
$ws.Range("E6").Value = "'5.23%"
$ws.Range("D7").Value = "'8.316"
$ws.Range("E7").Value = "'4.38%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9382"
$ws.Range("E8").Value = "'0.77%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1351"
$ws.Range("E9").Value = "'-9.59%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1977"
$ws.Range("E10").Value = "'2.48%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09071"
$ws.Range("E11").Value = "'0.23%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03493"
$ws.Range("E12").Value = "'-0.12%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09830"
$ws.Range("E13").Value = "'0.55%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001372"
$ws.Range("E14").Value = "'-2.33%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006322"
$ws.Range("E15").Value = "'5.97%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.689"
$ws.Range("E16").Value = "'-2.57%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.321"
$ws.Range("E17").Value = "'3.09%"
$ws.Range("E19").Value = "'2.05%"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("E20").Value = "'-0.02%"
$ws.Range("D21").Value = "'4.855"
$ws.Range("E21").Value = "'7.09%"
$ws.Range("D22").Value = "'0.2444"
$ws.Range("E22").Value = "'1.27%"
$ws.Range("D23").Value = "'0.04327"
$ws.Range("E23").Value = "'-1.00%"
$ws.Range("D24").Value = "'0.001225"
$ws.Range("E24").Value = "'-0.89%"
$ws.Range("D25").Value = "'0.004819"
$ws.Range("E25").Value = "'12.49%"
$ws.Range("E26").Value = "'-0.30%"
$ws.Range("D27").Value = "'0.0003990"
$ws.Range("E27").Value = "'-10.29%"
$ws.Range("D39").Value = "'0.02222"
$ws.Range("E39").Value = "'8.75%"
$ws.Range("D40").Value = "'0.05207"
$ws.Range("E40").Value = "'2.03%"
$ws.Range("D41").Value = "'0.007669"
$ws.Range("E41").Value = "'3.16%"
$ws.Range("D42").Value = "'0.009632"
$ws.Range("E42").Value = "'-6.43%"
$ws.Range("D43").Value = "'0.1404"
$ws.Range("E43").Value = "'4.01%"
$ws.Range("D44").Value = "'0.002087"
$ws.Range("E44").Value = "'-1.59%"
$ws.Range("D45").Value = "'0.008746"
$ws.Range("E45").Value = "'-3.97%"
$ws.Range("D46").Value = "'0.00006647"
$ws.Range("E46").Value = "'7.42%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.30%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.002880"
$ws.Range("E48").Value = "'-7.14%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.001686"
$ws.Range("E49").Value = "'5.26%"
$ws.Range("D50").Value = "'0.00002095"
$ws.Range("E50").Value = "'-0.30%"
$ws.Range("D51").Value = "'0.0001995"
$ws.Range("E51").Value = "'-0.30%"
